# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.404.09"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "1.838.74"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.62"
$ws.Range("E5").Value = "  -6.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5204"
$ws.Range("E7").Value = "  -1.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3257"
$ws.Range("E8").Value = "  -6.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06766"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.73"
$ws.Range("E10").Value = "  -7.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7651"
$ws.Range("E11").Value = "  -5.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07673"
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("D13").Value = "1.851.42"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.99"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.033"
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.99"
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007883"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").Value = "26.424.86"
$ws.Range("D21").Value = "2.067.81"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.564"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.460"
$ws.Range("E23").Value = "  -5.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.949"
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.257"
$ws.Range("E25").Value = "  -4.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.50"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.99"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.61"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.172"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.144"
$ws.Range("E31").Value = "  -4.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08707"
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04808"
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.126"
$ws.Range("E34").Value = "  -4.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.853"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6955"
$ws.Range("E36").Value = "  -5.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.069"
$ws.Range("E37").Value = "  -7.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01766"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.199"
$ws.Range("E39").Value = "  -8.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4868"
$ws.Range("E40").Value = "  -5.98%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8931"
$ws.Range("E41").Value = "  -7.30%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.75"
$ws.Range("E42").Value = "  -4.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.087"
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.712"
$ws.Range("E45").Value = "  -5.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4167"
$ws.Range("E46").Value = "  -8.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.031"
$ws.Range("E47").Value = "  -3.89%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05862"
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1244"
$ws.Range("E49").Value = "  -7.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.83"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8822"
$ws.Range("E51").Value = "  -0.64%  "
